$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (no longer present in the updated data)
$ws.Range("A8:T10").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.115543
$ws.Range("H2").Value = 0.346629
$ws.Range("I2").Value = 0.7111155332715143
$ws.Range("J2").Value = 0.7111155332715143
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05145833333333334
$ws.Range("N2").Value = 0.154375
$ws.Range("O2").Value = 0.2409462730781657
$ws.Range("P2").Value = 0.2409462730781657
$ws.Range("Q2").Value = 0.005945650208333333
$ws.Range("R2").Value = 0.05351085187500001
$ws.Range("S2").Value = 0.1713406374697637
$ws.Range("T2").Value = 0.1713406374697637

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.115543
$ws.Range("H3").Value = 0.346629
$ws.Range("I3").Value = 0.7111155332715143
$ws.Range("J3").Value = 0.7111155332715143
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1394176666666667
$ws.Range("N3").Value = 0.418253
$ws.Range("O3").Value = 0.6528032489312521
$ws.Range("P3").Value = 0.6528032489312521
$ws.Range("Q3").Value = 0.01610873545966667
$ws.Range("R3").Value = 0.144978619137
$ws.Range("S3").Value = 0.4642185304851245
$ws.Range("T3").Value = 0.4642185304851245

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl28"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.115543
$ws.Range("H4").Value = 0.346629
$ws.Range("I4").Value = 0.7111155332715143
$ws.Range("J4").Value = 0.7111155332715143
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02269166666666667
$ws.Range("N4").Value = 0.068075
$ws.Range("O4").Value = 0.1062504779905822
$ws.Range("P4").Value = 0.1062504779905822
$ws.Range("Q4").Value = 0.002621863241666667
$ws.Range("R4").Value = 0.023596769175
$ws.Range("S4").Value = 0.07555636531662617
$ws.Range("T4").Value = 0.07555636531662617

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl28"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04693833333333333
$ws.Range("H5").Value = 0.140815
$ws.Range("I5").Value = 0.2888844667284857
$ws.Range("J5").Value = 0.2888844667284857
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05145833333333334
$ws.Range("N5").Value = 0.154375
$ws.Range("O5").Value = 0.2409462730781657
$ws.Range("P5").Value = 0.2409462730781657
$ws.Range("Q5").Value = 0.002415368402777778
$ws.Range("R5").Value = 0.021738315625
$ws.Range("S5").Value = 0.069605635608402
$ws.Range("T5").Value = 0.069605635608402

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ccl28"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04693833333333333
$ws.Range("H6").Value = 0.140815
$ws.Range("I6").Value = 0.2888844667284857
$ws.Range("J6").Value = 0.2888844667284857
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1394176666666667
$ws.Range("N6").Value = 0.418253
$ws.Range("O6").Value = 0.6528032489312521
$ws.Range("P6").Value = 0.6528032489312521
$ws.Range("Q6").Value = 0.006544032910555555
$ws.Range("R6").Value = 0.058896296195
$ws.Range("S6").Value = 0.1885847184461277
$ws.Range("T6").Value = 0.1885847184461277

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ccl28"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04693833333333333
$ws.Range("H7").Value = 0.140815
$ws.Range("I7").Value = 0.2888844667284857
$ws.Range("J7").Value = 0.2888844667284857
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02269166666666667
$ws.Range("N7").Value = 0.068075
$ws.Range("O7").Value = 0.1062504779905822
$ws.Range("P7").Value = 0.1062504779905822
$ws.Range("Q7").Value = 0.001065109013888889
$ws.Range("R7").Value = 0.009585981124999999
$ws.Range("S7").Value = 0.03069411267395605
$ws.Range("T7").Value = 0.03069411267395605
